$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 147.5
$ws.Range("I31").Value = 96.666664
$ws.Range("J31").Value = 300
$ws.Range("K31").Value = 289.999992
$ws.Range("L31").Value = 900
$ws.Range("M31").Value = -59.99999200000002
$ws.Range("N31").Value = -1360

$ws.Range("H132").Value = 3168766
$ws.Range("I132").Value = 661568.1
$ws.Range("J132").Value = 22223470
$ws.Range("K132").Value = 1984704.3
$ws.Range("L132").Value = 66670410
$ws.Range("M132").Value = -1982174.3
$ws.Range("N132").Value = -66675470

$ws.Range("H135").Value = 35716468
$ws.Range("I135").Value = 1085.875
$ws.Range("J135").Value = 83336980
$ws.Range("K135").Value = 9772.875
$ws.Range("L135").Value = 750032820
$ws.Range("M135").Value = -7237.875
$ws.Range("N135").Value = -750037890

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20484.889
$ws.Range("I2").Value = 22914.25
$ws.Range("J2").Value = 1050
$ws.Range("K2").Value = 22914.25
$ws.Range("L2").Value = 1050
$ws.Range("M2").Value = -22801.25
$ws.Range("N2").Value = -1276

$ws.Range("H37").Value = 10006.789
$ws.Range("J37").Value = 14212.9
$ws.Range("L37").Value = 14212.9
$ws.Range("N37").Value = -14758.9

$ws.Range("H74").Value = 25533454
$ws.Range("I74").Value = 28572126
$ws.Range("J74").Value = 16670664
$ws.Range("K74").Value = 28572126
$ws.Range("L74").Value = 16670664
$ws.Range("M74").Value = -28571252
$ws.Range("N74").Value = -16672412

$ws.Range("H77").Value = 25533454
$ws.Range("I77").Value = 28572126
$ws.Range("J77").Value = 16670664
$ws.Range("K77").Value = 142860630
$ws.Range("L77").Value = 83353320
$ws.Range("M77").Value = -142856262
$ws.Range("N77").Value = -83362056

$ws.Range("H116").Value = 20484.889
$ws.Range("I116").Value = 22914.25
$ws.Range("J116").Value = 1050
$ws.Range("K116").Value = 22914.25
$ws.Range("L116").Value = 1050
$ws.Range("M116").Value = -20620.25
$ws.Range("N116").Value = -5638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20484.889
$ws.Range("I3").Value = 22914.25
$ws.Range("J3").Value = 1050
$ws.Range("K3").Value = 22914.25
$ws.Range("L3").Value = 1050
$ws.Range("M3").Value = -22800.25
$ws.Range("N3").Value = -1278

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1854717.2
$ws.Range("I31").Value = 2526427.8
$ws.Range("J31").Value = 7513.6665
$ws.Range("K31").Value = 2526427.8
$ws.Range("L31").Value = 7513.6665
$ws.Range("M31").Value = -2526132.8
$ws.Range("N31").Value = -8103.6665

$ws.Range("H34").Value = 1854717.2
$ws.Range("I34").Value = 2526427.8
$ws.Range("J34").Value = 7513.6665
$ws.Range("K34").Value = 2526427.8
$ws.Range("L34").Value = 7513.6665
$ws.Range("M34").Value = -2526225.8
$ws.Range("N34").Value = -7917.6665

$ws.Range("H50").Value = 14452.667
$ws.Range("J50").Value = 14452.667
$ws.Range("L50").Value = 14452.667
$ws.Range("N50").Value = -15702.667

$ws.Range("H51").Value = 9172.286
$ws.Range("J51").Value = 9172.286
$ws.Range("L51").Value = 9172.286
$ws.Range("N51").Value = -10644.286

$ws.Range("H59").Value = 16000.667
$ws.Range("J59").Value = 16000.667
$ws.Range("L59").Value = 16000.667
$ws.Range("N59").Value = -18290.667

$ws.Range("H60").Value = 7516.143
$ws.Range("J60").Value = 8385.5
$ws.Range("L60").Value = 8385.5
$ws.Range("N60").Value = -9407.5

$ws.Range("H61").Value = 9172.286
$ws.Range("J61").Value = 9172.286
$ws.Range("L61").Value = 9172.286
$ws.Range("N61").Value = -9868.286

$ws.Range("H68").Value = 18156.555
$ws.Range("J68").Value = 18156.555
$ws.Range("L68").Value = 18156.555
$ws.Range("N68").Value = -19654.555

$ws.Range("H70").Value = 18000
$ws.Range("J70").Value = 18000
$ws.Range("L70").Value = 18000
$ws.Range("N70").Value = -18630

$ws.Range("H71").Value = 18156.555
$ws.Range("J71").Value = 18156.555
$ws.Range("L71").Value = 54469.665
$ws.Range("N71").Value = -61957.665

$ws.Range("H73").Value = 18000
$ws.Range("J73").Value = 18000
$ws.Range("L73").Value = 18000
$ws.Range("N73").Value = -20184

$ws.Range("H74").Value = 20139.5
$ws.Range("J74").Value = 20139.5
$ws.Range("L74").Value = 20139.5
$ws.Range("N74").Value = -21887.5

$ws.Range("H77").Value = 20139.5
$ws.Range("J77").Value = 20139.5
$ws.Range("L77").Value = 60418.5
$ws.Range("N77").Value = -69154.5

$ws.Range("H99").Value = 14499.353
$ws.Range("I99").Value = 26330
$ws.Range("J99").Value = 11964.214
$ws.Range("K99").Value = 26330
$ws.Range("L99").Value = 11964.214
$ws.Range("M99").Value = -24832
$ws.Range("N99").Value = -14960.214

$ws.Range("H126").Value = 14499.353
$ws.Range("I126").Value = 26330
$ws.Range("J126").Value = 11964.214
$ws.Range("K126").Value = 78990
$ws.Range("L126").Value = 35892.642
$ws.Range("M126").Value = -76520
$ws.Range("N126").Value = -40832.642

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 808.75
$ws.Range("I13").Value = 468.33334
$ws.Range("J13").Value = 922.2222
$ws.Range("K13").Value = 1405.00002
$ws.Range("L13").Value = 2766.6666
$ws.Range("M13").Value = -1237.00002
$ws.Range("N13").Value = -3102.6666

$ws.Range("H68").Value = 1330.4
$ws.Range("I68").Value = 550.6667
$ws.Range("J68").Value = 2500
$ws.Range("K68").Value = 1652.0001
$ws.Range("L68").Value = 7500
$ws.Range("M68").Value = -841.0001
$ws.Range("N68").Value = -9122

$ws.Range("H71").Value = 1330.4
$ws.Range("I71").Value = 550.6667
$ws.Range("J71").Value = 2500
$ws.Range("K71").Value = 4956.0003
$ws.Range("L71").Value = 22500
$ws.Range("M71").Value = -900.0002999999997
$ws.Range("N71").Value = -30612

$ws.Range("H114").Value = 525.36365
$ws.Range("I114").Value = 258
$ws.Range("J114").Value = 846.2
$ws.Range("K114").Value = 774
$ws.Range("L114").Value = 2538.6
$ws.Range("M114").Value = 2480
$ws.Range("N114").Value = -9046.6

$ws.Range("H131").Value = 5061604
$ws.Range("I131").Value = 50100350
$ws.Range("J131").Value = 1070.191
$ws.Range("K131").Value = 150301050
$ws.Range("L131").Value = 3210.573
$ws.Range("M131").Value = -150296010
$ws.Range("N131").Value = -13290.573

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 9151095
$ws.Range("I132").Value = 8255363
$ws.Range("J132").Value = 12989949
$ws.Range("K132").Value = 24766089
$ws.Range("L132").Value = 38969847
$ws.Range("M132").Value = -24763559
$ws.Range("N132").Value = -38974907

$ws.Range("H136").Value = 19265.033
$ws.Range("J136").Value = 19265.033
$ws.Range("L136").Value = 57795.099
$ws.Range("N136").Value = -62895.099

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 29414698
$ws.Range("I22").Value = 860
$ws.Range("J22").Value = 35717664
$ws.Range("K22").Value = 860
$ws.Range("L22").Value = 35717664
$ws.Range("M22").Value = -565
$ws.Range("N22").Value = -35718254

$ws.Range("H27").Value = 29414698
$ws.Range("I27").Value = 860
$ws.Range("J27").Value = 35717664
$ws.Range("K27").Value = 860
$ws.Range("L27").Value = 35717664
$ws.Range("M27").Value = -753
$ws.Range("N27").Value = -35717878

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 290.26666
$ws.Range("I113").Value = 227.23077
$ws.Range("J113").Value = 700
$ws.Range("K113").Value = 681.69231
$ws.Range("L113").Value = 2100
$ws.Range("M113").Value = 1488.30769
$ws.Range("N113").Value = -6440
